# "Generate Report for Handoff"
# Replaces the two tracked source files (md1/md2) with a new pair, flips the
# status text to "Ready for handoff", refreshes the handoff timestamps, and
# drops the (now-empty) "Latest Target File"/"Latest Handback File" columns
# from the per-locale sheets.

$wb = $excel.ActiveWorkbook

$oldMd1 = "6ec5984e-4e80-45ba-b9c8-02d397fcd2e4"
$oldMd2 = "93dc6f00-5870-4bc0-9c25-05ab848c3ed4"
$newMd1 = "d90d6a7c-5ffc-4db8-9e07-497e5621c334"
$newMd2 = "ffffb0782bef-8311-4f94-aaa6-686e82cc04b3"

$newStatus = "Ready for handoff"
$overviewDate = "2016-41-11 12:41:56"

$newZhXlf = "$newMd1.8c699e8997456c7ca9e599d4d73f26a578a91859.zh-cn.xlf"
$newDeXlf = "$newMd1.8c699e8997456c7ca9e599d4d73f26a578a91859.de-de.xlf"
$zhHandoffDate = "2016-03-11 12:41:52"
$deHandoffDate = "2016-03-11 12:41:56"
$epoch = "0001-01-01 00:00:00"

function Set-HyperlinkText($sheet, $addr, $text) {
    foreach ($hl in $sheet.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
            return
        }
    }
}

function Remove-HyperlinkAt($sheet, $addr) {
    foreach ($hl in $sheet.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
            return
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newMd1.md"
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("D2").Value = $overviewDate

$ov.Range("A3").Value = "$newMd2.md"
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus
$ov.Range("D3").Value = $overviewDate

Set-HyperlinkText $ov '$A$2' "$newMd1.md"
Set-HyperlinkText $ov '$A$3' "$newMd2.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "$newMd1.md"
$zh.Range("C2").Value = $newStatus
$zh.Range("D2").Value = $newZhXlf
$zh.Range("E2").Value = $zhHandoffDate
$zh.Range("H2").Value = $epoch

$zh.Range("A3").Value = "$newMd2.md"
$zh.Range("C3").Value = $newStatus
$zh.Range("D3").Value = $newZhXlf
$zh.Range("E3").Value = $zhHandoffDate
$zh.Range("H3").Value = $epoch

Set-HyperlinkText $zh '$A$2' "$newMd1.md"
Set-HyperlinkText $zh '$D$2' $newZhXlf
Set-HyperlinkText $zh '$A$3' "$newMd2.md"
Set-HyperlinkText $zh '$D$3' $newZhXlf

Remove-HyperlinkAt $zh '$F$2'
Remove-HyperlinkAt $zh '$G$2'
Remove-HyperlinkAt $zh '$F$3'
Remove-HyperlinkAt $zh '$G$3'

$zh.Range("F2:G3").Clear()

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "$newMd1.md"
$de.Range("C2").Value = $newStatus
$de.Range("D2").Value = $newDeXlf
$de.Range("E2").Value = $deHandoffDate
$de.Range("H2").Value = $epoch

$de.Range("A3").Value = "$newMd2.md"
$de.Range("C3").Value = $newStatus
$de.Range("D3").Value = $newDeXlf
$de.Range("E3").Value = $deHandoffDate
$de.Range("H3").Value = $epoch

Set-HyperlinkText $de '$A$2' "$newMd1.md"
Set-HyperlinkText $de '$D$2' $newDeXlf
Set-HyperlinkText $de '$A$3' "$newMd2.md"
Set-HyperlinkText $de '$D$3' $newDeXlf

Remove-HyperlinkAt $de '$F$2'
Remove-HyperlinkAt $de '$G$2'
Remove-HyperlinkAt $de '$F$3'
Remove-HyperlinkAt $de '$G$3'

$de.Range("F2:G3").Clear()
